# Edit script: applies the diff to before.xlsx
#  1. Add new worksheet "엔씨소프트" at the end (sheetId=4, rId4)
#  2. Populate header row and 100 data rows (dates + remn_amt) on the new sheet
#  3. Update B102 ("remn_amt" for 2025-10-0x date row) on the three existing
#     sheets (카카오, NAVER, 농심) from 0 to the new observed values

$wb = $excel.ActiveWorkbook

# --- existing sheets: update the previously-zero B102 values ---
$wsKakao = $wb.Worksheets.Item("카카오")
$wsKakao.Range("B102").Value = 907245

$wsNaver = $wb.Worksheets.Item("NAVER")
$wsNaver.Range("B102").Value = 1218569

$wsNongshim = $wb.Worksheets.Item("농심")
$wsNongshim.Range("B102").Value = 133299

# --- add the new worksheet at the end of the workbook ---
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$newSheet = $wb.Worksheets.Add($null, $lastSheet)
$newSheet.Name = "엔씨소프트"

# Match the page margins used by the other sheets (0.75"/0.75"/1"/1"/0.5"/0.5")
$newSheet.PageSetup.LeftMargin = 54
$newSheet.PageSetup.RightMargin = 54
$newSheet.PageSetup.TopMargin = 72
$newSheet.PageSetup.BottomMargin = 72
$newSheet.PageSetup.HeaderMargin = 36
$newSheet.PageSetup.FooterMargin = 36

# Copy formatting (styles) from the first sheet so the new sheet reuses the
# same header style (bold/border/centered) and date-column number format
# instead of creating brand-new style entries.
$wsKakao.Range("A1:B1").Copy()
$newSheet.Range("A1:B1").PasteSpecial(-4122)

$wsKakao.Range("A2:B2").Copy()
$newSheet.Range("A2:B101").PasteSpecial(-4122)

# Header values
$newSheet.Range("A1").Value = "date"
$newSheet.Range("B1").Value = "remn_amt"

# Bulk-populate the 100 data rows (A2:B101) in one shot
$data = New-Object 'object[,]' 100,2

$data[0,0] = 45813
$data[0,1] = 238435
$data[1,0] = 45817
$data[1,1] = 222317
$data[2,0] = 45818
$data[2,1] = 230747
$data[3,0] = 45819
$data[3,1] = 246193
$data[4,0] = 45820
$data[4,1] = 239688
$data[5,0] = 45821
$data[5,1] = 228116
$data[6,0] = 45824
$data[6,1] = 207376
$data[7,0] = 45825
$data[7,1] = 190968
$data[8,0] = 45826
$data[8,1] = 203281
$data[9,0] = 45827
$data[9,1] = 202253
$data[10,0] = 45828
$data[10,1] = 202386
$data[11,0] = 45831
$data[11,1] = 209947
$data[12,0] = 45832
$data[12,1] = 229712
$data[13,0] = 45833
$data[13,1] = 224680
$data[14,0] = 45834
$data[14,1] = 218573
$data[15,0] = 45835
$data[15,1] = 215967
$data[16,0] = 45838
$data[16,1] = 241113
$data[17,0] = 45839
$data[17,1] = 238839
$data[18,0] = 45840
$data[18,1] = 231088
$data[19,0] = 45841
$data[19,1] = 231925
$data[20,0] = 45842
$data[20,1] = 228365
$data[21,0] = 45845
$data[21,1] = 225516
$data[22,0] = 45846
$data[22,1] = 229737
$data[23,0] = 45847
$data[23,1] = 221487
$data[24,0] = 45848
$data[24,1] = 217107
$data[25,0] = 45849
$data[25,1] = 211999
$data[26,0] = 45852
$data[26,1] = 207021
$data[27,0] = 45853
$data[27,1] = 203436
$data[28,0] = 45854
$data[28,1] = 199907
$data[29,0] = 45855
$data[29,1] = 189508
$data[30,0] = 45856
$data[30,1] = 187853
$data[31,0] = 45859
$data[31,1] = 187857
$data[32,0] = 45860
$data[32,1] = 186974
$data[33,0] = 45861
$data[33,1] = 180848
$data[34,0] = 45862
$data[34,1] = 175173
$data[35,0] = 45863
$data[35,1] = 176543
$data[36,0] = 45866
$data[36,1] = 177515
$data[37,0] = 45867
$data[37,1] = 176417
$data[38,0] = 45868
$data[38,1] = 176869
$data[39,0] = 45869
$data[39,1] = 176244
$data[40,0] = 45870
$data[40,1] = 172367
$data[41,0] = 45873
$data[41,1] = 175285
$data[42,0] = 45874
$data[42,1] = 170676
$data[43,0] = 45875
$data[43,1] = 170379
$data[44,0] = 45876
$data[44,1] = 176947
$data[45,0] = 45877
$data[45,1] = 169810
$data[46,0] = 45880
$data[46,1] = 170676
$data[47,0] = 45881
$data[47,1] = 190667
$data[48,0] = 45882
$data[48,1] = 182984
$data[49,0] = 45883
$data[49,1] = 180445
$data[50,0] = 45887
$data[50,1] = 177598
$data[51,0] = 45888
$data[51,1] = 177443
$data[52,0] = 45889
$data[52,1] = 155446
$data[53,0] = 45890
$data[53,1] = 151034
$data[54,0] = 45891
$data[54,1] = 150635
$data[55,0] = 45894
$data[55,1] = 148620
$data[56,0] = 45895
$data[56,1] = 153811
$data[57,0] = 45896
$data[57,1] = 153866
$data[58,0] = 45897
$data[58,1] = 156356
$data[59,0] = 45898
$data[59,1] = 152070
$data[60,0] = 45901
$data[60,1] = 143037
$data[61,0] = 45902
$data[61,1] = 134727
$data[62,0] = 45903
$data[62,1] = 135908
$data[63,0] = 45904
$data[63,1] = 138895
$data[64,0] = 45905
$data[64,1] = 141594
$data[65,0] = 45908
$data[65,1] = 140875
$data[66,0] = 45909
$data[66,1] = 142650
$data[67,0] = 45910
$data[67,1] = 141613
$data[68,0] = 45911
$data[68,1] = 143263
$data[69,0] = 45912
$data[69,1] = 157283
$data[70,0] = 45915
$data[70,1] = 165318
$data[71,0] = 45916
$data[71,1] = 156274
$data[72,0] = 45917
$data[72,1] = 159388
$data[73,0] = 45918
$data[73,1] = 155312
$data[74,0] = 45919
$data[74,1] = 150070
$data[75,0] = 45922
$data[75,1] = 149094
$data[76,0] = 45923
$data[76,1] = 153255
$data[77,0] = 45924
$data[77,1] = 154674
$data[78,0] = 45925
$data[78,1] = 148853
$data[79,0] = 45926
$data[79,1] = 144719
$data[80,0] = 45929
$data[80,1] = 152822
$data[81,0] = 45930
$data[81,1] = 148536
$data[82,0] = 45931
$data[82,1] = 150724
$data[83,0] = 45932
$data[83,1] = 157413
$data[84,0] = 45940
$data[84,1] = 155847
$data[85,0] = 45943
$data[85,1] = 153259
$data[86,0] = 45944
$data[86,1] = 154221
$data[87,0] = 45945
$data[87,1] = 154194
$data[88,0] = 45946
$data[88,1] = 160835
$data[89,0] = 45947
$data[89,1] = 153166
$data[90,0] = 45950
$data[90,1] = 157008
$data[91,0] = 45951
$data[91,1] = 153394
$data[92,0] = 45952
$data[92,1] = 145025
$data[93,0] = 45953
$data[93,1] = 142547
$data[94,0] = 45954
$data[94,1] = 142344
$data[95,0] = 45957
$data[95,1] = 149794
$data[96,0] = 45958
$data[96,1] = 153619
$data[97,0] = 45959
$data[97,1] = 156093
$data[98,0] = 45960
$data[98,1] = 157826
$data[99,0] = 45961
$data[99,1] = 156380

$newSheet.Range("A2:B101").Value = $data

Write-Output "Edit complete."
